# LOM3250 -> "Trabalho de Graduação II" update
#
# Applies the textual edits plus the new italic English paragraphs
# described by the commit diff.

$d = $word.ActiveDocument

function Insert-ItalicParagraphAfter($paragraph, [string]$text) {
    # Create a new, empty paragraph right after $paragraph, then type
    # $text into it and italicise only the run (not the paragraph mark).
    $paragraph.Range.InsertParagraphAfter() | Out-Null
    $newPara = $paragraph.Next()
    $start = $newPara.Range.Start
    $newPara.Range.InsertAfter($text)
    $end = $newPara.Range.End
    $textRange = $d.Range($start, $end - 1)
    $textRange.Font.Italic = $true
}

# 1) Title heading: "LOM3250 -  Trabalho de Graduação" -> "... II"
$d.Content.Find.Execute(
    "LOM3250 -  Trabalho de Graduação", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "LOM3250 -  Trabalho de Graduação II", 2) | Out-Null

# 2) Subtitle heading: "Undergraduate Work" -> "Graduation Monograph II"
$d.Content.Find.Execute(
    "Undergraduate Work", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "Graduation Monograph II", 2) | Out-Null

# 3) Activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/2012", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "Ativação: 01/01/2023", 2) | Out-Null

# 4) New italic English paragraph after the "Objetivos" body paragraph
$objPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("O Trabalho de Graduação (TG) tem por objetivo")) {
        $objPara = $p
        break
    }
}
Insert-ItalicParagraphAfter $objPara "The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer."

# 5) Docente(s) Responsável(eis) list: replace the single professor with two
$docentePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("519033 - Carlos Yujiro Shigue")) {
        $docentePara = $p
        break
    }
}
$docRange = $docentePara.Range
$docRange.MoveEnd(1, -1) | Out-Null
$docRange.Text = "5840730 - Antonio Jefferson da Silva Machado"
$docRange.Collapse(0)
$docRange.InsertAfter([char]11)
$docRange.Collapse(0)
$docRange.InsertAfter("1176388 - Luiz Tadeu Fernandes Eleno")

# 6) New italic English paragraph after "Programa resumido" body paragraph
$resumoPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Elaborar uma monografia de Trabalho de Graduação")) {
        $resumoPara = $p
        break
    }
}
Insert-ItalicParagraphAfter $resumoPara "Prepare a monograph of Undergraduate Work under the guidance of a professor and present it to a panel of examiners."

# 7) New italic English paragraph after "Programa" body paragraph
$programaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("O programa da disciplina será constituído")) {
        $programaPara = $p
        break
    }
}
Insert-ItalicParagraphAfter $programaPara "The course program will consist of the following steps: 1) Preparation and writing of a monograph on a previously defined and approved subject in the Undergraduate Work I discipline. 2) Definition and disclosure of the presentation date after delivery of the monograph in advance of at least , 15 working days. 3) Definition of the panel of examiners, consisting of the supervisor and at least two invited professionals, with training in engineering or related areas. 4) Presentation and evaluation of the TG. 5) Publication of the evaluation. In case of approval, the final copy of the monograph (printed and electronic copy) must be delivered with the agreement of the supervisor."

# 8) Requisitos: swap the prerequisite course
$d.Content.Find.Execute(
    "LOM3238 -  Projeto Integrado I  (Requisito)", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "LOM3267 -  Trabalho de Graduação I  (Requisito)", 2) | Out-Null
